# Apply the upstream "Add files via upload" edit to Dataset.xlsx / Hoja1.
#
# The change re-runs geocoding for a handful of "Lima" / "Ayacucho" rows,
# nudging latitude (col H) / longitude (col I) to more precise coordinates.
# A few of the updated latitude cells became live formulas (a plain unary
# negation of the literal) instead of static numbers; the rest stay as
# plain numeric literals. Column F is hidden, and the sheet view's
# selection / scroll position is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Latitude (H) / Longitude (I) updates -----------------------------
# Rows whose H value in the diff carries an <f> element are written via
# .Formula (a unary-negation formula); everything else is a plain .Value.

$ws.Range("H2").Value = -12.045005977682299
$ws.Range("I2").Value = -77.029964136886207

$ws.Range("H3").Value = -13.160994057531299
$ws.Range("I3").Value = -74.2325736178272

$ws.Range("H4").Value = -12.0578260465577
$ws.Range("I4").Value = -77.0406736466811

$ws.Range("H5").Value = -12.8268585509094
$ws.Range("I5").Value = -74.111507049831403

$ws.Range("H7").Formula = "=-12.0450059776823"
$ws.Range("I7").Value = -77.029964136886207

$ws.Range("H8").Value = -12.096712767937699
$ws.Range("I8").Value = -77.021037741032004

$ws.Range("H9").Value = -12.0463216510866
$ws.Range("I9").Value = -77.028377719694802

$ws.Range("H10").Formula = "=-12.0450059776823"
$ws.Range("I10").Value = -77.029964136886207

$ws.Range("H11").Formula = "=-12.0479451976752"
$ws.Range("I11").Value = -77.025359597612507

$ws.Range("H12").Value = -12.1186696082605
$ws.Range("I12").Value = -77.000673939473003

$ws.Range("H13").Formula = "=-12.0479451976752"
$ws.Range("I13").Value = -77.025359597612507

$ws.Range("H15").Value = -12.299061880493101
$ws.Range("I15").Value = -76.854659896952199

$ws.Range("H16").Formula = "=-12.0479451976752"
$ws.Range("I16").Value = -77.025359597612507

$ws.Range("H18").Formula = "=-12.0450059776823"
$ws.Range("I18").Value = -77.029964136886207

$ws.Range("H19").Formula = "=-12.0479451976752"
$ws.Range("I19").Value = -77.025359597612507

$ws.Range("H22").Formula = "=-12.0517187883729"
$ws.Range("I22").Value = -77.034671572863004

$ws.Range("H23").Formula = "=-12.0450059776823"
$ws.Range("I23").Value = -77.029964136886207

$ws.Range("H28").Formula = "=-12.0450059776823"
$ws.Range("I28").Value = -77.029964136886207

$ws.Range("H30").Formula = "=-12.0479451976752"
$ws.Range("I30").Value = -77.025359597612507

$ws.Range("H32").Formula = "=-12.0479451976752"
$ws.Range("I32").Value = -77.025359597612507

$ws.Range("H34").Formula = "=-12.0450059776823"
$ws.Range("I34").Value = -77.029964136886207

$ws.Range("H45").Formula = "=-12.0517187883729"
$ws.Range("I45").Value = -77.034671572863004

$ws.Range("H46").Formula = "=-12.0517187883729"
$ws.Range("I46").Value = -77.034671572863004

# --- Column F becomes hidden --------------------------------------------
$ws.Range("F1").EntireColumn.Hidden = $true

# --- Sheet view: scroll position + selection ----------------------------
# topLeftCell moves from A10 to C1, and the active selection moves from
# K23 to H17.
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("H17").Select()
